$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.00000000000047
$ws.Range("H2").Value = [double]"1.265211424074252e-16"
$ws.Range("I2").Value = 0.8681181715780015
$ws.Range("K2").Value = 41.14151352281606
$ws.Range("L2").Value = "[34.61512082869976, 47.66790621693235]"
$ws.Range("O2").Value = 1.62897396852804
$ws.Range("P2").Value = "[1.452868674633116, 1.8050792624229643]"
$ws.Range("S2").Value = 50.47748036621294
$ws.Range("T2").Value = "[46.19110162265448, 54.76385910977141]"
$ws.Range("W2").Value = 18.51851851851887
$ws.Range("X2").Value = 17.81781781781815
$ws.Range("Y2").Value = 19.21921921921959

# Row 3 updates
$ws.Range("E3").Value = 24.33000000000036
$ws.Range("H3").Value = [double]"1.265211424074252e-16"
$ws.Range("K3").Value = 45.36690039303943
$ws.Range("L3").Value = "[36.78683212919346, 53.946968656885396]"
$ws.Range("O3").Value = 0.748447499053424
$ws.Range("P3").Value = "[0.5597632555945768, 0.9371317425122712]"
$ws.Range("Q3").Value = [double]"3.512745649913995e-13"
$ws.Range("R3").Value = [double]"3.512745649913995e-13"
$ws.Range("S3").Value = 54.82039780873326
$ws.Range("T3").Value = "[50.41379912093982, 59.2269964965267]"
$ws.Range("W3").Value = 21.43183183183215
$ws.Range("X3").Value = 20.70120120120151
$ws.Range("Y3").Value = 22.1624624624628
